# Revert "install redux and create Redux store":
# The commit being reverted had appended a second screenshot (plus a
# blank separator paragraph) below the original architecture diagram.
# Undo that: drop the last inline picture's paragraph and the blank
# paragraph that precedes it, leaving the first diagram and the
# trailing bookmark paragraph untouched.

$d = $word.ActiveDocument

if ($d.InlineShapes.Count -gt 1) {
    $lastShape = $d.InlineShapes.Item($d.InlineShapes.Count)
    $shapeStart = $lastShape.Range.Start

    # Locate the paragraph that owns this inline shape.
    $paraIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Start -eq $shapeStart) {
            $paraIndex = $i
        }
    }

    if ($paraIndex -gt 0) {
        # Remove the picture's own paragraph.
        $d.Paragraphs.Item($paraIndex).Range.Delete()

        # Remove the now-orphaned blank paragraph that separated the
        # two pictures (if any).
        if ($paraIndex -gt 1) {
            $prev = $d.Paragraphs.Item($paraIndex - 1)
            if ($prev.Range.Text -eq "`r" -and $prev.Range.InlineShapes.Count -eq 0) {
                $prev.Range.Delete()
            }
        }
    }
}
